$d = $word.ActiveDocument

# 1. Scatter-mode bullet -> Pink ghost bullet
$d.Content.Find.Execute(
    "Now change ghost to move according to scatter mode. Research behavior and make sure it reacts accordingly.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pink ghost (middle) should move out of box and move in a similar fashion as red ghost for now.", 2)

# 2. Chase-mode bullet -> Blue/Orange up-down bullet
$d.Content.Find.Execute(
    "Enact chase mode and make sure ghost reacts accordingly.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Blue and Orange ghost should move up and down.", 2)

# 3. Cycle scatter/chase bullet -> Blue ghost center-of-box bullet
$d.Content.Find.Execute(
    "Cycle through scatter and chase mode. Follow guide lines.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Blue ghost should move to center of box and then up and out. Make this ghost move like the others.", 2)

# 4. Incorporate pink ghost -> Orange ghost center-of-box bullet
$d.Content.Find.Execute(
    "Incorporate pink ghost.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Orange ghost is the last to go and should also move to center of box then up and out.", 2)

# 5. Incorporate yellow ghost -> player/ghost collision detection
$d.Content.Find.Execute(
    "Incorporate yellow ghost.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Add player and ghost collision detection.", 2)

# 6. "Incorporate " + "Blue ghost." (two runs, one paragraph) -> single run
#    "Add death state and have Pac-Man animate death." -- Find/Replace merges
#    the two runs into one clean run and keeps the trailing bookmark
#    (_GoBack) tracking the edit point.
$d.Content.Find.Execute(
    "Incorporate Blue ghost.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Add death state and have Pac-Man animate death.", 2)

# Now append the new bullet's text right at the bookmark (end of the text
# we just wrote) so the bookmark tracks forward naturally, then split the
# paragraph into two at that same point using a Find/Replace that inserts
# a paragraph mark ("^p"). This reproduces the exact structure Word itself
# produces (new list paragraph, bookmark at the tail of the new paragraph).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Range.InsertAfter("Have level restart. Check with video or play game to see how level resets.")

$d.Content.Find.Execute(
    "death.Have level restart",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "death.^pHave level restart", 2)
